$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): Right count 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): Right count 78 -> 130
$ws.Range("B12").Value = 130

# Update the "Corr/total" marks text 76/84 -> 130/140
$ws.Range("E12").Value = "130/140"
